# Remove the three duplicate/erroneous order rows for E/S PUNTA HERMOSA
# (rows 19-21) so that each delivery compartment maps to a single customer.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A19:A21").EntireRow.Delete()
